# Update sample tag data and update 'create' function in Tag model
#
# 1. Tag sheet: dedupe the tag list (case-insensitive dedupe of "Academic",
#    "Research", lowercase "life"/"place" typos), trim from 20 rows to 14.
# 2. PollTag sheet: remap the tag_id foreign key column to the
#    deduplicated tag ids.
# 3. View state: PollTag becomes the active/selected tab (was Answer).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Tag sheet: rewrite the (deduplicated) tag names, then drop the now
# redundant trailing rows.
# ---------------------------------------------------------------------
$tagSheet = $wb.Worksheets.Item("Tag")

$tagNames = @(
    '"Life"',
    '"Food"',
    '"Place"',
    '"Academic"',
    '"Information Systems"',
    '"Professor"',
    '"Course"',
    '"Career"',
    '"Time"',
    '"Research"',
    '"Discrimination"',
    '"Favorite"',
    '"Freshman"',
    '"Recommendation"'
)

for ($i = 0; $i -lt $tagNames.Length; $i++) {
    $row = $i + 2
    $tagSheet.Range("B" + $row).Value = $tagNames[$i]
}

# Old sheet had 20 tag rows (rows 2-21); the new deduped list only needs
# rows 2-15, so remove the leftover rows 16-21.
$tagSheet.Rows("16:21").Delete()
$tagSheet.Range("A1:B15").Select()

# ---------------------------------------------------------------------
# PollTag sheet: remap tag_id (column C) from the old tag ids to the new
# deduplicated tag ids.
# ---------------------------------------------------------------------
$pollTagSheet = $wb.Worksheets.Item("PollTag")

$newTagIds = @{
    8  = 2
    9  = 6
    10 = 3
    11 = 3
    12 = 7
    13 = 8
    14 = 9
    15 = 9
    16 = 10
    17 = 0
    18 = 11
    19 = 12
    20 = 0
    21 = 13
}

foreach ($row in $newTagIds.Keys) {
    $pollTagSheet.Range("C" + $row).Value = $newTagIds[$row]
}

# ---------------------------------------------------------------------
# View state: PollTag is now the selected/active sheet (previously it
# was Answer).
# ---------------------------------------------------------------------
$pollTagSheet.Activate()
$pollTagSheet.Range("A1:C21").Select()
